$d = $word.ActiveDocument

$d.Content.Find.Execute("811×9=7299", $true, $false, $false, $false, $false, $true, 1, $false, "675×7=4725", 2)
$d.Content.Find.Execute("576×8=4608", $true, $false, $false, $false, $false, $true, 1, $false, "219×4=876", 2)
$d.Content.Find.Execute("411×9=3699", $true, $false, $false, $false, $false, $true, 1, $false, "778×9=7002", 2)
$d.Content.Find.Execute("159×2=318", $true, $false, $false, $false, $false, $true, 1, $false, "778×4=3112", 2)
$d.Content.Find.Execute("515×7=3605", $true, $false, $false, $false, $false, $true, 1, $false, "749×6=4494", 2)
$d.Content.Find.Execute("379×7=2653", $true, $false, $false, $false, $false, $true, 1, $false, "848×4=3392", 2)
$d.Content.Find.Execute("411×8=3288", $true, $false, $false, $false, $false, $true, 1, $false, "429×7=3003", 2)
$d.Content.Find.Execute("275×8=2200", $true, $false, $false, $false, $false, $true, 1, $false, "647×2=1294", 2)
$d.Content.Find.Execute("778×6=4668", $true, $false, $false, $false, $false, $true, 1, $false, "189×2=378", 2)
$d.Content.Find.Execute("157×5=785", $true, $false, $false, $false, $false, $true, 1, $false, "155×6=930", 2)
$d.Content.Find.Execute("573×9=5157", $true, $false, $false, $false, $false, $true, 1, $false, "737×6=4422", 2)
$d.Content.Find.Execute("806×6=4836", $true, $false, $false, $false, $false, $true, 1, $false, "970×3=2910", 2)
$d.Content.Find.Execute("987×5=4935", $true, $false, $false, $false, $false, $true, 1, $false, "981×9=8829", 2)
$d.Content.Find.Execute("591×7=4137", $true, $false, $false, $false, $false, $true, 1, $false, "942×6=5652", 2)
$d.Content.Find.Execute("142×4=568", $true, $false, $false, $false, $false, $true, 1, $false, "965×5=4825", 2)
$d.Content.Find.Execute("434×4=1736", $true, $false, $false, $false, $false, $true, 1, $false, "843×7=5901", 2)
$d.Content.Find.Execute("850×7=5950", $true, $false, $false, $false, $false, $true, 1, $false, "402×9=3618", 2)
$d.Content.Find.Execute("149×7=1043", $true, $false, $false, $false, $false, $true, 1, $false, "417×7=2919", 2)
$d.Content.Find.Execute("165×8=1320", $true, $false, $false, $false, $false, $true, 1, $false, "242×2=484", 2)
$d.Content.Find.Execute("212×9=1908", $true, $false, $false, $false, $false, $true, 1, $false, "169×3=507", 2)
$d.Content.Find.Execute("331×7=2317", $true, $false, $false, $false, $false, $true, 1, $false, "547×6=3282", 2)
$d.Content.Find.Execute("864×5=4320", $true, $false, $false, $false, $false, $true, 1, $false, "949×3=2847", 2)
$d.Content.Find.Execute("616×9=5544", $true, $false, $false, $false, $false, $true, 1, $false, "137×2=274", 2)
$d.Content.Find.Execute("345×5=1725", $true, $false, $false, $false, $false, $true, 1, $false, "628×3=1884", 2)
$d.Content.Find.Execute("943×9=8487", $true, $false, $false, $false, $false, $true, 1, $false, "309×7=2163", 2)
